$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New SOFR OIS data rows for 2025-09-08 through 2025-09-12
$dates = @("2025-09-08","2025-09-08","2025-09-08","2025-09-08","2025-09-08","2025-09-08","2025-09-08","2025-09-09","2025-09-09","2025-09-09","2025-09-09","2025-09-09","2025-09-09","2025-09-09","2025-09-10","2025-09-10","2025-09-10","2025-09-10","2025-09-10","2025-09-10","2025-09-10","2025-09-11","2025-09-11","2025-09-11","2025-09-11","2025-09-11","2025-09-11","2025-09-11","2025-09-12","2025-09-12","2025-09-12","2025-09-12","2025-09-12","2025-09-12","2025-09-12")
$tenors = @("1Y","2Y","3Y","5Y","10Y","20Y","30Y","1Y","2Y","3Y","5Y","10Y","20Y","30Y","1Y","2Y","3Y","5Y","10Y","20Y","30Y","1Y","2Y","3Y","5Y","10Y","20Y","30Y","1Y","2Y","3Y","5Y","10Y","20Y","30Y")
$rates = @(0.035889,0.032494,0.031652,0.03203,0.035202,0.038809,0.038736,0.036396,0.032997,0.032083,0.032355,0.03548,0.039063,0.038958,0.036191,0.032851,0.031948,0.032167,0.035115,0.038665,0.038564,0.035877,0.032752,0.031925,0.032132,0.034901,0.038387,0.038277,0.036111,0.033047,0.032276,0.032528,0.03532,0.038721,0.038535)

$startRow = 317
$endRow = $startRow + $dates.Length - 1

# Force column A to text format first so the "YYYY-MM-DD" strings are
# stored as plain text (matching the source data) instead of being
# auto-converted into date serial numbers by Excel.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $tenors[$i]
    $ws.Cells.Item($r, 3).Value = $rates[$i]
}

# Restore default (no explicit) cell formatting so the new rows use the
# same unstyled cells as the rest of the data table.
$ws.Range("A$startRow`:A$endRow").ClearFormats()

Write-Host "Added $($dates.Length) rows starting at row $startRow"
